# registration form and readme update
#
# Re-creates the user edits captured in the commit:
#   - "To Do" sheet: scrolled the frozen view down a bit and left the
#     cursor on B14 (was C12); also moved the two stray "x" marker
#     cells from column D to column F on rows 21-22 so they line up
#     with every other row's marker column.
#   - "Rename", "Sheet2" and "Kylie's Advice" sheets: opened Page Setup
#     and confirmed Portrait orientation (the "To Do" sheet already had
#     this set).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "To Do" sheet
# ---------------------------------------------------------------------
$todo = $wb.Worksheets.Item("To Do")

# Two cells with the "x" complete-marker were sitting in column D on
# rows 21 and 22 while every other row uses column F - slide them over.
$marker21 = $todo.Range("D21").Formula
$todo.Range("F21").Formula = $marker21
$todo.Range("D21").ClearContents()

$marker22 = $todo.Range("D22").Formula
$todo.Range("F22").Formula = $marker22
$todo.Range("D22").ClearContents()

# Scroll the frozen (header-locked) view down a few rows and leave the
# selection on B14.
$todo.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$todo.Range("B14").Select()

# ---------------------------------------------------------------------
# Page setup - set Portrait orientation on the other three sheets
# ---------------------------------------------------------------------
$xlPortrait = 1
$pageSetupSheets = @("Rename", "Sheet2", "Kylie's Advice")
foreach ($sheetName in $pageSetupSheets) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $sheet.PageSetup.Orientation = $xlPortrait
}

# Leave the originally active sheet/selection in place.
$todo.Activate()
$todo.Range("B14").Select()
